$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at X (24th column), pushing the existing FRAUD_RISK
# column (and any columns to its right) one position to the right.
$ws.Columns.Item(24).Insert()

# New column header
$ws.Cells.Item(1, 24).Value = "WEEKEND_APPR_PROCESS_START"

# Fill the new column's data rows (2-6) with "N/A"
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 24).Value = "N/A"
}
